# userDB.xlsx: rename the sheet, drop the "role" column from both the
# worksheet and its backing table, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "工作表1" -> "userDB"
$ws.Name = "userDB"

# Remove the third table column ("role"); this also deletes the
# underlying worksheet column (C) and its data/ shared strings, and
# shrinks the table ref from A1:C7 to A1:B7.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item(3).Delete()

# Move the saved selection to H13.
[void]$ws.Range("H13").Select()
